$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly price feed: each week contributes a "Primera"/"Segunda"
# pair of rows. A new week (date 44665) is being inserted right after the
# existing row 45, which pushes every following week down by two rows. The
# two rows that fall off the bottom (old rows 136/137) become new rows 138/139.

$lastCol = 18   # column R

# 1) Append new rows 138 and 139 using the CURRENT (pre-shift) contents of
#    rows 136 and 137, before anything below row 48 gets overwritten.
for ($col = 1; $col -le $lastCol; $col++) {
    $ws.Cells.Item(138, $col).Value = $ws.Cells.Item(136, $col).Value2
    $ws.Cells.Item(139, $col).Value = $ws.Cells.Item(137, $col).Value2
}
$ws.Cells.Item(138, 4).NumberFormat = $ws.Cells.Item(136, 4).NumberFormat
$ws.Cells.Item(139, 4).NumberFormat = $ws.Cells.Item(137, 4).NumberFormat

# 2) Shift every week from the old row 46 downward by two rows: new[r] = old[r-2].
#    Walk from the bottom up so we never overwrite a row before reading it.
for ($r = 137; $r -ge 48; $r--) {
    for ($col = 1; $col -le $lastCol; $col++) {
        $ws.Cells.Item($r, $col).Value = $ws.Cells.Item($r - 2, $col).Value2
    }
}

# 3) Rows 46/47 keep their original data but represent the newly-added week.
$ws.Cells.Item(46, 4).Value = 44665
$ws.Cells.Item(47, 4).Value = 44665
